# Chile "primera-division" 2023 results sheet update.
#
# 1) Three pairs of adjacent fixtures (rows 26/27, 42/43, 45/46) that share
#    the same kickoff date in column E were re-ordered by the scraper --
#    their F:V odds/result data is swapped while the leading index/meta
#    columns (A:E) stay untouched.
# 2) Two newly-scraped fixtures are appended as rows 209/210.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap F:V between row 26 and row 27 (reorder of same-date fixtures)
$ws.Range("F26").Value = "Union La Calera"
$ws.Range("F27").Value = "Copiapo"
$ws.Range("G26").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H26").Value = "Curico Unido"
$ws.Range("H27").Value = "Palestino"
$ws.Range("I26").Value = 2
$ws.Range("I27").Value = 3
$ws.Range("J26").Value = 2.34
$ws.Range("J27").Value = 3.43
$ws.Range("K26").Value = "04/02/2023 22:12"
$ws.Range("K27").Value = "06/02/2023 00:42"
$ws.Range("L26").Value = 2.04
$ws.Range("L27").Value = 3.18
$ws.Range("M26").Value = "11/02/2023 00:57"
$ws.Range("M27").Value = "11/02/2023 00:53"
$ws.Range("N26").Value = 3.62
$ws.Range("N27").Value = 3.91
$ws.Range("O26").Value = "04/02/2023 22:12"
$ws.Range("O27").Value = "06/02/2023 00:42"
$ws.Range("P26").Value = 3.42
$ws.Range("P27").Value = 3.62
$ws.Range("Q26").Value = "11/02/2023 00:58"
$ws.Range("Q27").Value = "11/02/2023 00:53"
$ws.Range("R26").Value = 3.03
$ws.Range("R27").Value = 1.94
$ws.Range("S26").Value = "04/02/2023 22:12"
$ws.Range("S27").Value = "06/02/2023 00:42"
$ws.Range("T26").Value = 3.99
$ws.Range("T27").Value = 2.26
$ws.Range("U26").Value = "11/02/2023 00:57"
$ws.Range("U27").Value = "11/02/2023 00:53"
$ws.Range("V26").Value = "https://www.betexplorer.com/football/chile/primera-division/union-la-calera-curico-unido/4Mz9Ngnb/"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/chile/primera-division/copiapo-palestino/WGWP0Wvo/"

# Swap F:V between row 42 and row 43 (reorder of same-date fixtures)
$ws.Range("F42").Value = "Union La Calera"
$ws.Range("F43").Value = "Magallanes"
$ws.Range("G42").Value = 1
$ws.Range("G43").Value = 2
$ws.Range("H42").Value = "Nublense"
$ws.Range("H43").Value = "Cobresal"
$ws.Range("I42").Value = 1
$ws.Range("I43").Value = 1
$ws.Range("J42").Value = 2.4
$ws.Range("J43").Value = 2.74
$ws.Range("K42").Value = "20/02/2023 22:12"
$ws.Range("K43").Value = "20/02/2023 16:42"
$ws.Range("L42").Value = 2
$ws.Range("L43").Value = 2.67
$ws.Range("M42").Value = "26/02/2023 21:57"
$ws.Range("M43").Value = "26/02/2023 21:57"
$ws.Range("N42").Value = 3.28
$ws.Range("N43").Value = 3.54
$ws.Range("O42").Value = "20/02/2023 22:12"
$ws.Range("O43").Value = "20/02/2023 16:42"
$ws.Range("P42").Value = 3.6
$ws.Range("P43").Value = 3.55
$ws.Range("Q42").Value = "26/02/2023 21:57"
$ws.Range("Q43").Value = "26/02/2023 21:57"
$ws.Range("R42").Value = 2.95
$ws.Range("R43").Value = 2.59
$ws.Range("S42").Value = "20/02/2023 22:12"
$ws.Range("S43").Value = "20/02/2023 16:42"
$ws.Range("T42").Value = 3.89
$ws.Range("T43").Value = 2.66
$ws.Range("U42").Value = "26/02/2023 21:53"
$ws.Range("U43").Value = "26/02/2023 21:58"
$ws.Range("V42").Value = "https://www.betexplorer.com/football/chile/primera-division/union-la-calera-nublense/rXo8WlVS/"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/chile/primera-division/magallanes-cobresal/IP0rxg1i/"

# Swap F:V between row 45 and row 46 (reorder of same-date fixtures)
$ws.Range("F45").Value = "Everton"
$ws.Range("F46").Value = "Copiapo"
$ws.Range("G45").Value = 3
$ws.Range("G46").Value = 2
$ws.Range("H45").Value = "A. Italiano"
$ws.Range("H46").Value = "O'Higgins"
$ws.Range("I45").Value = 1
$ws.Range("I46").Value = 2
$ws.Range("J45").Value = 1.87
$ws.Range("J46").Value = 3.21
$ws.Range("K45").Value = "20/02/2023 00:41"
$ws.Range("K46").Value = "20/02/2023 16:42"
$ws.Range("L45").Value = 1.76
$ws.Range("L46").Value = 2.51
$ws.Range("M45").Value = "28/02/2023 00:28"
$ws.Range("M46").Value = "28/02/2023 00:05"
$ws.Range("N45").Value = 3.62
$ws.Range("N46").Value = 3.48
$ws.Range("O45").Value = "20/02/2023 00:41"
$ws.Range("O46").Value = "20/02/2023 16:42"
$ws.Range("P45").Value = 3.73
$ws.Range("P46").Value = 3.32
$ws.Range("Q45").Value = "28/02/2023 00:29"
$ws.Range("Q46").Value = "28/02/2023 00:11"
$ws.Range("R45").Value = 4.48
$ws.Range("R46").Value = 2.16
$ws.Range("S45").Value = "20/02/2023 00:41"
$ws.Range("S46").Value = "20/02/2023 16:42"
$ws.Range("T45").Value = 4.99
$ws.Range("T46").Value = 3
$ws.Range("U45").Value = "28/02/2023 00:29"
$ws.Range("U46").Value = "28/02/2023 00:11"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/chile/primera-division/everton-a-italiano/EifizXV3/"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/chile/primera-division/copiapo-o-higgins/08d3XUFM/"

# Append new rows 209 and 210 (new fixtures), copying row 208's formatting first
$ws.Range("A208:V208").Copy() | Out-Null
$ws.Range("A209:V210").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column D ("temporada") holds numeric-looking text ("2023") like every other
# row in the sheet -- force Text format first so Excel keeps it as a string
# instead of coercing it to the number 2023.
$ws.Range("D209:D210").NumberFormat = "@"

# Row 209
$ws.Range("A209").Value = 208
$ws.Range("B209").Value = "chile"
$ws.Range("C209").Value = "primera-division"
$ws.Range("D209").Value = "2023"
$ws.Range("E209").Value = 45241.9375
$ws.Range("F209").Value = "Cobresal"
$ws.Range("G209").Value = 3
$ws.Range("H209").Value = "A. Italiano"
$ws.Range("I209").Value = 4
$ws.Range("J209").Value = 1.51
$ws.Range("K209").Value = "03/11/2023 22:42"
$ws.Range("L209").Value = 1.68
$ws.Range("M209").Value = "11/11/2023 22:07"
$ws.Range("N209").Value = 4.38
$ws.Range("O209").Value = "03/11/2023 22:42"
$ws.Range("P209").Value = 4.12
$ws.Range("Q209").Value = "11/11/2023 22:07"
$ws.Range("R209").Value = 6.41
$ws.Range("S209").Value = "03/11/2023 22:42"
$ws.Range("T209").Value = 5.01
$ws.Range("U209").Value = "11/11/2023 22:16"
$ws.Range("V209").Value = "https://www.betexplorer.com/football/chile/primera-division/cobresal-a-italiano/nL8BEppO/"

# Row 210
$ws.Range("A210").Value = 209
$ws.Range("B210").Value = "chile"
$ws.Range("C210").Value = "primera-division"
$ws.Range("D210").Value = "2023"
$ws.Range("E210").Value = 45242.04166666666
$ws.Range("F210").Value = "Nublense"
$ws.Range("G210").Value = 1
$ws.Range("H210").Value = "Palestino"
$ws.Range("I210").Value = 1
$ws.Range("J210").Value = 2.67
$ws.Range("K210").Value = "04/11/2023 01:43"
$ws.Range("L210").Value = 2.85
$ws.Range("M210").Value = "12/11/2023 00:51"
$ws.Range("N210").Value = 3.45
$ws.Range("O210").Value = "04/11/2023 01:43"
$ws.Range("P210").Value = 3.19
$ws.Range("Q210").Value = "12/11/2023 00:52"
$ws.Range("R210").Value = 2.65
$ws.Range("S210").Value = "04/11/2023 01:43"
$ws.Range("T210").Value = 2.71
$ws.Range("U210").Value = "12/11/2023 00:51"
$ws.Range("V210").Value = "https://www.betexplorer.com/football/chile/primera-division/nublense-palestino/QB7FD4aU/"
